$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.828.63'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +7.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.761.44'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.59%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +1.97%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.82'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.91%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9963'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.71%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3834'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3646'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +5.39%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.12'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +17.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.231'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07653'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +5.46%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9994'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.78'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.467'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +7.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.124'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +5.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.768.73'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +6.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001160'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.80%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9968'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.94%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06850'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '87.35'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +7.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.67'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.500'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.75'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +5.89%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '25.813.56'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.426'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.20%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.960'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +8.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.76'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '154.84'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.54%  '

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.963.95'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +6.61%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.05'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +6.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.206'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +20.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.093'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +14.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.297'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.46%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '14.06'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +13.97%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.817'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.52%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08718'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.668'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06761'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +6.60%  '

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.344'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.26%  '

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02474'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2227'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.301'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.90%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6581'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.23'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +7.92%  '

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9952'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.71%  '

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6365'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.73%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.913'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.10%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.170'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +8.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.67'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.95%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07503'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.64'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +6.18%  '
